$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-11 Tuesday", "2025-03-12 Wednesday"),
    @("919×9=8271", "683×6=4098"),
    @("999×9=8991", "376×3=1128"),
    @("937×6=5622", "513×2=1026"),
    @("686×2=1372", "782×7=5474"),
    @("727×9=6543", "801×3=2403"),
    @("371×8=2968", "250×2=500"),
    @("597×2=1194", "714×4=2856"),
    @("429×7=3003", "833×6=4998"),
    @("686×6=4116", "872×2=1744"),
    @("325×2=650", "968×2=1936"),
    @("417×7=2919", "450×5=2250"),
    @("254×5=1270", "811×5=4055"),
    @("803×9=7227", "285×4=1140"),
    @("973×2=1946", "197×5=985"),
    @("441×6=2646", "405×2=810"),
    @("654×6=3924", "790×5=3950"),
    @("102×5=510", "344×7=2408"),
    @("136×2=272", "506×8=4048"),
    @("245×2=490", "750×2=1500"),
    @("272×3=816", "719×3=2157"),
    @("753×2=1506", "140×2=280"),
    @("373×5=1865", "883×7=6181"),
    @("649×7=4543", "373×4=1492"),
    @("195×2=390", "281×7=1967"),
    @("935×2=1870", "492×8=3936")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced '$old' -> '$new': $found"
}
